# Natmi following Dr Hou advice
#
# The "M2" cluster is now included among the sending/target clusters for
# the Il4 -> Il13ra1 ligand-receptor pair, so the result table grows from
# 3x4 (ECs/FAPs/sCs x ECs/FAPs/M2/sCs) to a full 4x4 grid over
# ECs/FAPs/M2/sCs x ECs/FAPs/M2/sCs, and every downstream
# detection-rate / expression / specificity value is recomputed
# accordingly (12 data rows -> 16 data rows, A2:T17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Il4","Il13ra1","ECs",3.0,1.0,0.9334919999999999,2.800476,0.209095668602663,0.209095668602663,3.0,1.0,22.2939,66.8817,0.3209125068732572,0.3209125068732573,20.81117729879999,187.3005956892,0.06710141518762042,0.06710141518762042),
    @("ECs","Il4","Il13ra1","FAPs",3.0,1.0,0.9334919999999999,2.800476,0.209095668602663,0.209095668602663,3.0,1.0,22.81307866666667,68.439236,0.3283858932002398,0.3283858932002398,21.295826430704,191.662437876336,0.06866406789838682,0.06866406789838682),
    @("ECs","Il4","Il13ra1","M2",3.0,1.0,0.9334919999999999,2.800476,0.209095668602663,0.209095668602663,3.0,1.0,19.143989,57.43196699999999,0.2755706942950341,0.2755706942950342,17.87076057958799,160.836845216292,0.05762063857092022,0.05762063857092022),
    @("ECs","Il4","Il13ra1","sCs",3.0,1.0,0.9334919999999999,2.800476,0.209095668602663,0.209095668602663,3.0,1.0,5.219369333333334,15.658108,0.07513090563146878,0.07513090563146878,4.872239517712,43.850155659408,0.01570954694573554,0.01570954694573554),
    @("FAPs","Il4","Il13ra1","ECs",3.0,1.0,1.964334,5.893002,0.4399970552387631,0.4399970552387631,3.0,1.0,22.2939,66.8817,0.3209125068732572,0.3209125068732573,43.7926657626,394.1339918634,0.1412005580135225,0.1412005580135225),
    @("FAPs","Il4","Il13ra1","FAPs",3.0,1.0,1.964334,5.893002,0.4399970552387631,0.4399970552387631,3.0,1.0,22.81307866666667,68.439236,0.3283858932002398,0.3283858932002398,44.812506069608,403.312554626472,0.1444888259900565,0.1444888259900565),
    @("FAPs","Il4","Il13ra1","M2",3.0,1.0,1.964334,5.893002,0.4399970552387631,0.4399970552387631,3.0,1.0,19.143989,57.43196699999999,0.2755706942950341,0.2755706942950342,37.60518848832599,338.446696394934,0.1212502939999164,0.1212502939999165),
    @("FAPs","Il4","Il13ra1","sCs",3.0,1.0,1.964334,5.893002,0.4399970552387631,0.4399970552387631,3.0,1.0,5.219369333333334,15.658108,0.07513090563146878,0.07513090563146878,10.252584640024,92.273261760216,0.03305737723526767,0.03305737723526767),
    @("M2","Il4","Il13ra1","ECs",2.0,0.6666666666666666,0.1531806666666667,0.459542,0.03431139625585257,0.03431139625585256,3.0,1.0,22.2939,66.8817,0.3209125068732572,0.3209125068732573,3.414994464599999,30.7349501814,0.01101095618678734,0.01101095618678734),
    @("M2","Il4","Il13ra1","FAPs",2.0,0.6666666666666666,0.1531806666666667,0.459542,0.03431139625585257,0.03431139625585256,3.0,1.0,22.81307866666667,68.439236,0.3283858932002398,0.3283858932002398,3.494522598879111,31.450703389912,0.01126737850642551,0.01126737850642551),
    @("M2","Il4","Il13ra1","M2",2.0,0.6666666666666666,0.1531806666666667,0.459542,0.03431139625585257,0.03431139625585256,3.0,1.0,19.143989,57.43196699999999,0.2755706942950341,0.2755706942950342,2.932488997679333,26.392400979114,0.009455215288457327,0.009455215288457327),
    @("M2","Il4","Il13ra1","sCs",2.0,0.6666666666666666,0.1531806666666667,0.459542,0.03431139625585257,0.03431139625585256,3.0,1.0,5.219369333333334,15.658108,0.07513090563146878,0.07513090563146878,0.7995064740595555,7.195558266536001,0.002577846274182391,0.00257784627418239),
    @("sCs","Il4","Il13ra1","ECs",3.0,1.0,1.413418666666667,4.240256,0.3165958799027214,0.3165958799027214,3.0,1.0,22.2939,66.8817,0.3209125068732572,0.3209125068732573,31.5106144128,283.5955297152,0.101599577485327,0.101599577485327),
    @("sCs","Il4","Il13ra1","FAPs",3.0,1.0,1.413418666666667,4.240256,0.3165958799027214,0.3165958799027214,3.0,1.0,22.81307866666667,68.439236,0.3283858932002398,0.3283858932002398,32.24443123160178,290.199881084416,0.103965620805371,0.103965620805371),
    @("sCs","Il4","Il13ra1","M2",3.0,1.0,1.413418666666667,4.240256,0.3165958799027214,0.3165958799027214,3.0,1.0,19.143989,57.43196699999999,0.2755706942950341,0.2755706942950342,27.05847140706133,243.526242663552,0.08724454643574019,0.08724454643574019),
    @("sCs","Il4","Il13ra1","sCs",3.0,1.0,1.413418666666667,4.240256,0.3165958799027214,0.3165958799027214,3.0,1.0,5.219369333333334,15.658108,0.07513090563146878,0.07513090563146878,7.37715404396089,66.39438639564801,0.02378613517628319,0.02378613517628318)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $colLetter = [char](65 + $c)
        $ws.Range("$colLetter$r").Value = $row[$c]
    }
}
